# Scheduled refresh of the cryptos price list (coinranking.com snapshot),
# mirroring the GitHub Actions bot commit: updates Price (D) and Volume(1h) (E)
# columns for the affected rows. Values are written as plain text, matching
# the inlineStr cell type already used throughout the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.203.66"
$ws.Range("E2").Value = "  +1.72%  "

$ws.Range("D3").Value = "2.449.17"
$ws.Range("E3").Value = "  +1.03%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.535"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.81%  "

$ws.Range("D9").Value = "2.451.58"
$ws.Range("E9").Value = "  +1.08%  "

$ws.Range("E10").Value = "  +2.40%  "

$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("E12").Value = "  +2.15%  "

$ws.Range("E13").Value = "  +1.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.85%  "

$ws.Range("E15").Value = "  +3.60%  "

$ws.Range("D17").Value = "63.226.32"
$ws.Range("E17").Value = "  +1.87%  "

$ws.Range("D18").Value = "2.444.66"
$ws.Range("E18").Value = "  -1.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.79%  "

$ws.Range("E22").Value = "  +1.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "622.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000103"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.33%  "

$ws.Range("D29").Value = "2.590.84"
$ws.Range("E29").Value = "  +1.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.80%  "

$ws.Range("E33").Value = "  -1.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.36%  "

$ws.Range("E35").Value = "  +7.41%  "

$ws.Range("E36").Value = "  -0.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("E38").Value = "  +0.72%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.80%  "

$ws.Range("E41").Value = "  +0.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "145.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +15.74%  "

$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.42%  "

$ws.Range("E46").Value = "  +2.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0541"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.22%  "

$ws.Range("E49").Value = "  +0.71%  "

$ws.Range("E50").Value = "  +3.24%  "

$ws.Range("E51").Value = "  +0.35%  "
